# Apply the "Tested a6300 LL ev0" edit:
#  - Sheets "23" and "34" each have a block of 8 "measurement" rows (2-9)
#    whose Left/Top (B/C) shrink by 1 and Right/Down (D/E) grow by 1
#    (tolerance widened by 1px). The H:K helper columns are formulas
#    (=B-$G$2 etc.) so they recompute automatically.
#  - The ORIGINAL (pre-edit) B:E values of rows 2-9 get archived as
#    literal numbers into the H:K columns of a later "Backup" block
#    (rows 21-28 on sheet "23", rows 24-31 on sheet "34").
#  - On sheet "23" only, the now-redundant H:K formulas on rows 10-18
#    are cleared.
#  - Selection / active-sheet/tab bookkeeping: sheet "23" selection
#    moves to D8 and loses tabSelected; sheet "34" becomes the active
#    tab with selection F8.

$wb = $excel.ActiveWorkbook

$sheetNames = @("23", "34")
$backupStartRow = @{ "23" = 21; "34" = 24 }
$backupStyle    = @{ "23" = $null; "34" = "1" }

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # 1) Snapshot the original B:E values for rows 2-9 before they change.
    $orig = @{}
    for ($r = 2; $r -le 9; $r++) {
        $orig[$r] = @(
            $ws.Cells.Item($r, 2).Value(),
            $ws.Cells.Item($r, 3).Value(),
            $ws.Cells.Item($r, 4).Value(),
            $ws.Cells.Item($r, 5).Value()
        )
    }

    # 2) Shrink Left/Top by 1, grow Right/Down by 1 (tolerance +1).
    for ($r = 2; $r -le 9; $r++) {
        $vals = $orig[$r]
        $ws.Cells.Item($r, 2).Value = $vals[0] - 1
        $ws.Cells.Item($r, 3).Value = $vals[1] - 1
        $ws.Cells.Item($r, 4).Value = $vals[2] + 1
        $ws.Cells.Item($r, 5).Value = $vals[3] + 1
    }

    # 3) Archive the pre-edit values as literals into the Backup block.
    $startRow = $backupStartRow[$name]
    $style = $backupStyle[$name]
    for ($i = 0; $i -le 7; $i++) {
        $srcRow = 2 + $i
        $dstRow = $startRow + $i
        $vals = $orig[$srcRow]
        $ws.Cells.Item($dstRow, 8).Value = $vals[0]
        $ws.Cells.Item($dstRow, 9).Value = $vals[1]
        $ws.Cells.Item($dstRow, 10).Value = $vals[2]
        $ws.Cells.Item($dstRow, 11).Value = $vals[3]
        if ($style) {
            $ws.Cells.Item($dstRow, 8).Style = $style
            $ws.Cells.Item($dstRow, 9).Style = $style
            $ws.Cells.Item($dstRow, 10).Style = $style
            $ws.Cells.Item($dstRow, 11).Style = $style
        }
    }
}

# 4) Sheet "23": the H:K helper formulas for rows 10-18 are no longer needed.
$ws23 = $wb.Worksheets.Item("23")
$ws23.Range("H10:K18").ClearContents()

# 5) Selection / active tab bookkeeping.
$ws23.Activate()
$ws23.Range("D8").Select()

$ws34 = $wb.Worksheets.Item("34")
$ws34.Activate()
$ws34.Range("F8").Select()
